$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '64.157.67'
$ws.Range("E2").Value = '  +0.38%  '
$ws.Range("D3").Value = '3.150.87'
$ws.Range("E3").Value = '  +0.45%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '''590.16'
$ws.Range("E5").Value = '  -0.25%  '
$ws.Range("D6").Value = '''147.53'
$ws.Range("E6").Value = '  +0.43%  '
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").Value = '3.147.97'
$ws.Range("E8").Value = '  +0.61%  '
$ws.Range("E9").Value = '  -1.28%  '
$ws.Range("E10").Value = '  -2.92%  '
$ws.Range("D11").Value = '''5.88'
$ws.Range("E11").Value = '  +2.78%  '
$ws.Range("D12").Value = '''0.461'
$ws.Range("E12").Value = '  -1.68%  '
$ws.Range("E13").Value = '  -2.96%  '
$ws.Range("D14").Value = '''37.23'
$ws.Range("E14").Value = '  +2.95%  '
$ws.Range("D15").Value = '3.670.33'
$ws.Range("E15").Value = '  +0.41%  '
$ws.Range("E16").Value = '  -1.24%  '
$ws.Range("D17").Value = '63.927.19'
$ws.Range("E17").Value = '  +0.12%  '
$ws.Range("D18").Value = '3.146.61'
$ws.Range("E18").Value = '  +0.39%  '
$ws.Range("D19").Value = '''7.15'
$ws.Range("E19").Value = '  -0.54%  '
$ws.Range("D20").Value = '''466.28'
$ws.Range("E20").Value = '  -0.04%  '
$ws.Range("D21").Value = '''14.37'
$ws.Range("E21").Value = '  +0.54%  '
$ws.Range("D22").Value = '''0.734'
$ws.Range("E22").Value = '  +0.19%  '
$ws.Range("D23").Value = '''7.44'
$ws.Range("E23").Value = '  -1.76%  '
$ws.Range("D24").Value = '''13.05'
$ws.Range("E24").Value = '  -2.03%  '
$ws.Range("D25").Value = '''2.34'
$ws.Range("E25").Value = '  +8.28%  '
$ws.Range("D26").Value = '''80.96'
$ws.Range("E26").Value = '  -1.93%  '
$ws.Range("E27").Value = '  +0.07%  '
$ws.Range("D28").Value = '''9.65'
$ws.Range("E28").Value = '  +10.77%  '
$ws.Range("D29").Value = '''2.70'
$ws.Range("E29").Value = '  -0.74%  '
$ws.Range("B30").Value = 'ImmutableX'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D30").Value = '''2.22'
$ws.Range("E30").Value = '  -0.33%  '
$ws.Range("B31").Value = 'FirstDigitalUSD'
$ws.Range("C31").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D31").Value = '''1.00'
$ws.Range("E31").Value = '  -0.05%  '
$ws.Range("D32").Value = '''7.21'
$ws.Range("E32").Value = '  +5.15%  '
$ws.Range("D33").Value = '''27.20'
$ws.Range("E33").Value = '  +0.39%  '
$ws.Range("E34").Value = '  +0.25%  '
$ws.Range("D35").Value = '0.0₃0855'
$ws.Range("E35").Value = '  -1.94%  '
$ws.Range("D36").Value = '''1.07'
$ws.Range("E36").Value = '  +1.33%  '
$ws.Range("E37").Value = '  -3.75%  '
$ws.Range("D38").Value = '''6.08'
$ws.Range("E38").Value = '  -0.98%  '
$ws.Range("D39").Value = '''3.29'
$ws.Range("E39").Value = '  -4.87%  '
$ws.Range("D40").Value = '''51.72'
$ws.Range("E40").Value = '  +1.70%  '
$ws.Range("D41").Value = '''448.36'
$ws.Range("E41").Value = '  -1.84%  '
$ws.Range("D42").Value = '''8.95'
$ws.Range("E42").Value = '  +2.48%  '
$ws.Range("D43").Value = '''0.291'
$ws.Range("E43").Value = '  +4.80%  '
$ws.Range("E44").Value = '  -0.15%  '
$ws.Range("D45").Value = '2.913.14'
$ws.Range("E45").Value = '  -0.31%  '
$ws.Range("D46").Value = '''39.65'
$ws.Range("E46").Value = '  +14.71%  '
$ws.Range("E47").Value = '  -2.54%  '
$ws.Range("D48").Value = '''125.87'
$ws.Range("E48").Value = '  -2.53%  '
$ws.Range("B50").Value = 'ThetaToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D50").Value = '''2.24'
$ws.Range("E50").Value = '  +2.20%  '
$ws.Range("B51").Value = 'Stellar'
$ws.Range("C51").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D51").Value = '''0.111'
$ws.Range("E51").Value = '  -0.73%  '
